# Update column F ("想去人数" / "people interested" counts) on the
# "展览" and "全部类型" worksheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# Row => new F-column value for the "展览" sheet.
$sheetExhibitUpdates = @{
    2  = 1142
    3  = 855
    4  = 281
    5  = 51
    8  = 2382
    9  = 7755
    10 = 923
    11 = 446
    12 = 385
    13 = 158
    14 = 430
    16 = 162
    17 = 7968
    18 = 319
    19 = 1385
    20 = 159
    23 = 172
    24 = 325
    25 = 169
    29 = 27
    30 = 425
    31 = 1157
    35 = 83
    37 = 81
    38 = 70
}

# Row => new F-column value for the "全部类型" sheet (a couple of values
# differ slightly from the "展览" sheet's copy of the same rows).
$sheetAllTypesUpdates = @{
    2  = 1142
    3  = 855
    4  = 281
    5  = 51
    8  = 2383
    9  = 7755
    10 = 923
    11 = 446
    12 = 385
    13 = 158
    14 = 430
    16 = 162
    17 = 7969
    18 = 319
    19 = 1385
    20 = 159
    23 = 172
    24 = 325
    25 = 169
    29 = 27
    30 = 425
    31 = 1157
    35 = 83
    37 = 81
    38 = 70
}

$wsExhibit = $wb.Worksheets.Item("展览")
foreach ($row in $sheetExhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $sheetExhibitUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheetAllTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $sheetAllTypesUpdates[$row]
}
